$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text: volume number and report week dates ---
$ws.Range("A8").Value = "Volume 32   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/15/2025  Through  12/21/2025"

# --- Helper: set a numeric cell's value while inheriting number format from a
#     reference cell that already carries the desired style (keeps existing
#     cellXfs entries instead of minting new ones). ---
function Set-NumCell {
    param(
        [string]$Address,
        [double]$Value,
        [string]$FormatFromAddress
    )
    $ws.Range($Address).NumberFormat = $ws.Range($FormatFromAddress).NumberFormat
    $ws.Range($Address).Value = $Value
}

# --- Weekly crime-stat table refresh (rows 15-30) ---
Set-NumCell "C15" 1 "G15"
Set-NumCell "D15" 2 "G15"
Set-NumCell "E15" -50 "H15"
Set-NumCell "F15" 1 "G15"
Set-NumCell "G15" 4 "G15"
Set-NumCell "H15" -75 "H15"
Set-NumCell "I15" 20 "G15"
Set-NumCell "J15" 18 "G15"
Set-NumCell "K15" 11.111111111111 "H15"
Set-NumCell "L15" 42.857142857142 "H15"
Set-NumCell "M15" -9.090909090909 "H15"
Set-NumCell "N15" -77.52808988764 "H15"
Set-NumCell "C16" 2 "G15"
Set-NumCell "D16" 6 "G15"
Set-NumCell "E16" -66.666666666666 "H15"
Set-NumCell "I16" 235 "G15"
Set-NumCell "J16" 225 "G15"
Set-NumCell "K16" 4.444444444444 "H15"
Set-NumCell "L16" 12.980769230769 "H15"
Set-NumCell "M16" -15.162454873646 "H15"
Set-NumCell "N16" -81.909160892994 "H15"
Set-NumCell "C17" 4 "G15"
Set-NumCell "D17" 8 "G15"
Set-NumCell "E17" -50 "H15"
Set-NumCell "F17" 26 "G15"
Set-NumCell "G17" 29 "G15"
Set-NumCell "H17" -10.344827586206 "H15"
Set-NumCell "I17" 311 "G15"
Set-NumCell "J17" 344 "G15"
Set-NumCell "K17" -9.593023255813 "H15"
Set-NumCell "L17" -7.715133531157 "H15"
Set-NumCell "M17" 42.009132420091 "H15"
Set-NumCell "N17" -62.484921592279 "H15"
Set-NumCell "C18" 2 "G15"
Set-NumCell "D18" 4 "G15"
Set-NumCell "E18" -50 "H15"
Set-NumCell "G18" 14 "G15"
Set-NumCell "H18" -35.714285714285 "H15"
Set-NumCell "I18" 140 "G15"
Set-NumCell "J18" 122 "G15"
Set-NumCell "K18" 14.754098360655 "H15"
Set-NumCell "L18" -5.405405405405 "H15"
Set-NumCell "M18" -22.651933701657 "H15"
Set-NumCell "N18" -92.635455023671 "H15"
Set-NumCell "D19" 12 "G15"
Set-NumCell "E19" -16.666666666666 "H15"
Set-NumCell "F19" 51 "G15"
Set-NumCell "H19" 8.510638297872 "H15"
Set-NumCell "I19" 572 "G15"
Set-NumCell "J19" 628 "G15"
Set-NumCell "K19" -8.917197452229 "H15"
Set-NumCell "L19" 3.249097472924 "H15"
Set-NumCell "M19" 50.923482849604 "H15"
Set-NumCell "N19" -50.859106529209 "H15"
Set-NumCell "D20" 4 "G15"
Set-NumCell "E20" -50 "H15"
Set-NumCell "F20" 7 "G15"
Set-NumCell "G20" 11 "G15"
Set-NumCell "H20" -36.363636363636 "H15"
Set-NumCell "I20" 193 "G15"
Set-NumCell "J20" 164 "G15"
Set-NumCell "K20" 17.682926829268 "H15"
Set-NumCell "L20" -13.063063063063 "H15"
Set-NumCell "M20" 114.444444444444 "H15"
Set-NumCell "N20" -89.253897550111 "H15"
Set-NumCell "C21" 21 "C21"
Set-NumCell "D21" 36 "C21"
Set-NumCell "E21" -41.666666666666 "E21"
Set-NumCell "F21" 107 "C21"
Set-NumCell "G21" 121 "C21"
Set-NumCell "H21" -11.570247933884 "E21"
Set-NumCell "I21" 1474 "C21"
Set-NumCell "J21" 1507 "C21"
Set-NumCell "K21" -2.189781021897 "E21"
Set-NumCell "L21" -0.874243443174 "E21"
Set-NumCell "M21" 25.767918088737 "E21"
Set-NumCell "N21" -79.320987654321 "E21"
Set-NumCell "F22" 2 "G15"
Set-NumCell "G22" 2 "G15"
Set-NumCell "H22" 0 "H15"
Set-NumCell "J22" 29 "G15"
Set-NumCell "K22" -34.482758620689 "H15"
Set-NumCell "M22" -26.923076923076 "H15"
Set-NumCell "D24" 18 "G15"
Set-NumCell "E24" 16.666666666666 "H15"
Set-NumCell "F24" 103 "G15"
Set-NumCell "G24" 110 "G15"
Set-NumCell "H24" -6.363636363636 "H15"
Set-NumCell "I24" 1354 "G15"
Set-NumCell "J24" 1433 "G15"
Set-NumCell "K24" -5.512909979064 "H15"
Set-NumCell "L24" -2.519798416126 "H15"
Set-NumCell "M24" 119.093851132686 "H15"
Set-NumCell "C25" 12 "G15"
Set-NumCell "D25" 7 "G15"
Set-NumCell "E25" 71.428571428571 "H15"
Set-NumCell "G25" 49 "G15"
Set-NumCell "H25" -18.367346938775 "H15"
Set-NumCell "I25" 579 "G15"
Set-NumCell "J25" 681 "G15"
Set-NumCell "K25" -14.977973568281 "H15"
Set-NumCell "L25" -8.24088748019 "H15"
Set-NumCell "C26" 9 "G15"
Set-NumCell "D26" 2 "G15"
Set-NumCell "E26" 350 "H15"
Set-NumCell "G26" 34 "G15"
Set-NumCell "H26" 29.411764705882 "H15"
Set-NumCell "I26" 566 "G15"
Set-NumCell "J26" 555 "G15"
Set-NumCell "K26" 1.981981981981 "H15"
Set-NumCell "L26" 9.477756286266 "H15"
Set-NumCell "M26" 5.597014925373 "H15"
Set-NumCell "C27" 1 "G15"
Set-NumCell "D27" 2 "G15"
Set-NumCell "E27" -50 "H15"
Set-NumCell "F27" 2 "G15"
Set-NumCell "G27" 4 "G15"
Set-NumCell "I27" 26 "G15"
Set-NumCell "J27" 28 "G15"
Set-NumCell "K27" -7.142857142857 "H15"
Set-NumCell "L27" 4 "H15"
Set-NumCell "C28" 1 "G15"
Set-NumCell "D28" 2 "G15"
Set-NumCell "E28" -50 "H15"
Set-NumCell "F28" 5 "G15"
Set-NumCell "G28" 4 "G15"
Set-NumCell "H28" 25 "H15"
Set-NumCell "I28" 53 "G15"
Set-NumCell "J28" 61 "G15"
Set-NumCell "K28" -13.11475409836 "H15"
Set-NumCell "L28" -5.357142857142 "H15"
Set-NumCell "D29" 1 "G15"
Set-NumCell "E29" -100 "H15"
Set-NumCell "G29" 1 "G15"
Set-NumCell "H29" 0 "H15"
Set-NumCell "J29" 21 "G15"
Set-NumCell "K29" -61.904761904761 "H15"
Set-NumCell "L29" -55.555555555555 "H15"
Set-NumCell "N29" -95.4802259887 "H15"
Set-NumCell "D30" 1 "G15"
Set-NumCell "E30" -100 "H15"
Set-NumCell "G30" 1 "G15"
Set-NumCell "H30" 0 "H15"
Set-NumCell "J30" 18 "G15"
Set-NumCell "K30" -55.555555555555 "H15"
Set-NumCell "L30" -46.666666666666 "H15"
Set-NumCell "N30" -94.736842105263 "H15"
